$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "*" value in column C for the rows that indicate a completed/checked
# grading element (collection view functionality for the user and database
# drop function).
$ws.Range("C6").Value = "*"
$ws.Range("C11").Value = "*"
$ws.Range("C12").Value = "*"
$ws.Range("C15").Value = "*"
$ws.Range("C21").Value = "*"
$ws.Range("C22").Value = "*"
$ws.Range("C32").Value = "*"

# Move the active selection to C33
$ws.Range("C33").Select()
